# Add a new job posting row (Job_Id = 14) to Sheet1, right after the last
# existing row (row 14), following the same column layout as the rest of
# the table: Job_Id, Jd_Title, Job_Description, Total_Years_Min_Exp,
# Total_Years_Max_Exp, Linked_Posted, Resume_received, Resume_downloaded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 15

$ws.Cells.Item($newRow, 1).Value = 14              # Job_Id
$ws.Cells.Item($newRow, 2).Value = "AI Developer"  # Jd_Title
$ws.Cells.Item($newRow, 3).Value = "sfef"          # Job_Description
$ws.Cells.Item($newRow, 4).Value = 1               # Total_Years_Min_Exp
$ws.Cells.Item($newRow, 5).Value = 4               # Total_Years_Max_Exp
$ws.Cells.Item($newRow, 6).Value = 0               # Linked_Posted
$ws.Cells.Item($newRow, 7).Value = 0               # Resume_received
$ws.Cells.Item($newRow, 8).Value = 0               # Resume_downloaded
